# Update "想去人数" (want-to-go count) figures on the "展览" and "全部类型"
# sheets to the freshly scraped values.

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 33
$wsExpo.Range("F4").Value = 16135
$wsExpo.Range("F7").Value = 729
$wsExpo.Range("F8").Value = 15542
$wsExpo.Range("F10").Value = 9171
$wsExpo.Range("F11").Value = 447
$wsExpo.Range("F13").Value = 1025
$wsExpo.Range("F19").Value = 78
$wsExpo.Range("F20").Value = 588
$wsExpo.Range("F24").Value = 1135
$wsExpo.Range("F27").Value = 31
$wsExpo.Range("F33").Value = 2
$wsExpo.Range("F36").Value = 342
$wsExpo.Range("F39").Value = 5643
$wsExpo.Range("F40").Value = 5240

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 33
$wsAll.Range("F4").Value = 16136
$wsAll.Range("F7").Value = 729
$wsAll.Range("F8").Value = 15542
$wsAll.Range("F10").Value = 9171
$wsAll.Range("F11").Value = 447
$wsAll.Range("F13").Value = 1025
$wsAll.Range("F19").Value = 78
$wsAll.Range("F20").Value = 588
$wsAll.Range("F24").Value = 1135
$wsAll.Range("F27").Value = 31
$wsAll.Range("F35").Value = 2
$wsAll.Range("F38").Value = 342
$wsAll.Range("F41").Value = 5643
$wsAll.Range("F43").Value = 5240
